$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Торговля" (Trade) column already lists shop area / food-service seats (E4/E5);
# extend it with two more retail-turnover indicators, for all municipalities.
# Copy the look of the existing data cell (E5: centered text, shaded fill, border)
# onto the two new cells before filling them in.
$ws.Range("E5").Copy()
$ws.Range("E6:E7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("E6").Value = "Обор. роз. (кроме авто.) - retailturnover (тыс. руб.) (id8201003)"
$ws.Range("E7").Value = "Обор. Общепит - foodservturnover (тыс. руб.) (id8201006)"

# Widen column E so the longer labels are readable.
$ws.Columns("E").ColumnWidth = 65.5

# Leave the cursor where it ended up after typing the new entries.
$ws.Range("E8").Select()
